# UX e CSS menor - cabeçalhos melhorados
# Atualiza os dados de telemetria (duração, tempo total, velocidades,
# localização/coordenadas) refletindo a nova extração do relatório.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "8:19:06"
$ws.Range("H3").Value = "16 days 2:25:52"
$ws.Range("K3").Value = "98 km/h"
$ws.Range("G4").Value = "1:08:38"
$ws.Range("H4").Value = "15 days 23:57:35"
$ws.Range("K4").Value = "104 km/h"
$ws.Range("G6").Value = "0:02:06"
$ws.Range("H6").Value = "5 days 23:22:38"
$ws.Range("K6").Value = "71 km/h"
$ws.Range("G8").Value = "7:15:59"
$ws.Range("H8").Value = "29 days 1:34:17"
$ws.Range("K8").Value = "109 km/h"
$ws.Range("G9").Value = "17:30:43"
$ws.Range("H9").Value = "26 days 11:24:20"
$ws.Range("C10").Value = "22.07.2025 11:05:52"
$ws.Range("D10").Value = "55B, Subsetor Leste-6, Ribeirão Preto, SP, Brazil"
$ws.Range("E10").Value = -21.2048116
$ws.Range("F10").Value = -47.757315
$ws.Range("G10").Value = "3:50:01"
$ws.Range("H10").Value = "22 days 4:28:56"
$ws.Range("J10").Value = "75 km/h"
$ws.Range("L10").Value = "40 km/h"
$ws.Range("G11").Value = "5:57:58"
$ws.Range("H11").Value = "29 days 5:49:33"
$ws.Range("K11").Value = "115 km/h"
$ws.Range("C13").Value = "22.07.2025 08:51:55"
$ws.Range("D13").Value = "Sp-328, Ribeirão Preto, SP, Brazil"
$ws.Range("E13").Value = -21.0958533
$ws.Range("F13").Value = -47.7963533
$ws.Range("G13").Value = "14:06:15"
$ws.Range("H13").Value = "27 days 12:17:06"
$ws.Range("J13").Value = "102 km/h"
$ws.Range("G14").Value = "4:46:32"
$ws.Range("H14").Value = "9 days 23:44:22"
$ws.Range("G15").Value = "8:09:02"
$ws.Range("H15").Value = "13 days 7:06:23"
$ws.Range("G16").Value = "13:11:34"
$ws.Range("H16").Value = "22 days 2:32:33"
$ws.Range("C18").Value = "22.07.2025 10:06:50"
$ws.Range("D18").Value = "Avenida Gen. Euclydes De Figueiredo, Subsetor Norte-10, Ribeirão Preto, SP 14070-270, Brazil"
$ws.Range("E18").Value = -21.1112716
$ws.Range("F18").Value = -47.790145
$ws.Range("G18").Value = "14:36:40"
$ws.Range("H18").Value = "29 days 6:21:20"
$ws.Range("I18").Value = "{'t': '140 km/h', 'y': -20.4904, 'x': -42.1855633, 'u': 401879415}"
$ws.Range("J18").Value = "54 km/h"
$ws.Range("K18").Value = "93 km/h"
$ws.Range("G19").Value = "14:49:57"
$ws.Range("H19").Value = "14 days 8:57:10"
$ws.Range("K19").Value = "106 km/h"
$ws.Range("G20").Value = "7:03:02"
$ws.Range("H20").Value = "26 days 5:35:52"
$ws.Range("K20").Value = "105 km/h"
$ws.Range("G21").Value = "8:27:11"
$ws.Range("H21").Value = "14 days 6:37:00"
$ws.Range("G22").Value = "8:15:24"
$ws.Range("H22").Value = "27 days 5:50:48"
$ws.Range("K22").Value = "108 km/h"
$ws.Range("G23").Value = "9:00:27"
$ws.Range("H23").Value = "27 days 2:16:17"
$ws.Range("K23").Value = "106 km/h"
$ws.Range("G24").Value = "3:43:47"
$ws.Range("H24").Value = "14 days 0:23:52"
$ws.Range("I24").Value = "{'t': '151 km/h', 'y': -22.6260566, 'x': -50.5993383, 'u': 401929767}"
$ws.Range("K24").Value = "120 km/h"
$ws.Range("G25").Value = "11:11:10"
$ws.Range("H25").Value = "28 days 1:39:56"
$ws.Range("K25").Value = "101 km/h"
$ws.Range("G26").Value = "1:55:20"
$ws.Range("H26").Value = "15 days 7:40:16"
$ws.Range("G27").Value = "5:43:01"
$ws.Range("H27").Value = "26 days 22:59:33"
$ws.Range("G28").Value = "12:54:19"
$ws.Range("H28").Value = "26 days 2:25:13"
$ws.Range("I28").Value = "{'t': '189 km/h', 'y': -21.7579216, 'x': -48.0730216, 'u': 401897329}"
$ws.Range("K28").Value = "123 km/h"
$ws.Range("G29").Value = "6:47:54"
$ws.Range("H29").Value = "28 days 2:13:18"
$ws.Range("K29").Value = "112 km/h"
$ws.Range("G30").Value = "14:47:22"
$ws.Range("H30").Value = "20 days 19:51:41"
$ws.Range("K30").Value = "93 km/h"
$ws.Range("G33").Value = "0:28:18"
$ws.Range("H33").Value = "5 days 5:02:23"
$ws.Range("K33").Value = "88 km/h"
$ws.Range("G34").Value = "5:22:00"
$ws.Range("H34").Value = "27 days 0:57:28"
$ws.Range("K34").Value = "113 km/h"
$ws.Range("G35").Value = "16:29:00"
$ws.Range("H35").Value = "27 days 2:37:47"
$ws.Range("K35").Value = "105 km/h"
$ws.Range("G36").Value = "15:11:47"
$ws.Range("H36").Value = "15 days 13:48:18"
$ws.Range("K36").Value = "104 km/h"
